# Fix: adjust behavior based on old commit
#
# Inserts a new "Sheet1" worksheet right before "DxRenderer" documenting a
# new _outputSharedHandle / OutputSharedHandle plumbing change across
# DxRenderer.cs / IImageProcessor.cs / NativeImageProcessor.cs, and makes
# it the active sheet (mirroring the tabSelected move away from
# ImageProcCudaDLL).

$wb = $excel.ActiveWorkbook

$dxRenderer = $wb.Worksheets.Item("DxRenderer")

# New sheet goes immediately before "DxRenderer".
$ws = $wb.Worksheets.Add($dxRenderer)
$ws.Name = "Sheet1"

# --- populate content -------------------------------------------------
# Row 4 : new private field
$ws.Range("C4").Value = "_outputSharedHandle"
$ws.Range("B4").Value = "IntPtr"

# Row 5 : new public property + native-side equivalents
$ws.Range("C5").Value = "OutputSharedHandle"
$ws.Range("A5").Value = "public"

# Row 3 : existing field shown for reference
$ws.Range("A3").Value = "private readonly"
$ws.Range("B3").Value = "Texture2D"
$ws.Range("C3").Value = "_outputSharedTex"

# Row 4 finish
$ws.Range("A4").Value = "private readonly"

# Row 2 : file/class headers
$ws.Range("A2").Value = "DxRenderer.cs"
$ws.Range("D2").Value = "IImageProcessor.cs"

# Row 5 continue
$ws.Range("D5").Value = "System.IntPtr"
$ws.Range("E5").Value = "outputDxSharedHandle"
$ws.Range("G5").Value = "outputDxSharedHandle"

# Row 2 finish
$ws.Range("F2").Value = "NativeImageProcessor.cs"
$ws.Range("H2").Value = "NativeImageProc"

# Row 5 continue (native decl)
$ws.Range("I5").Value = "outSharedHandle"
$ws.Range("K5").Value = "outSharedHandle"
$ws.Range("J5").Value = "void*"

# Row 5 remaining IntPtr fills
$ws.Range("B5").Value = "IntPtr"
$ws.Range("F5").Value = "IntPtr"
$ws.Range("H5").Value = "IntPtr"

# Row 6
$ws.Range("A6").Value = "OpenSharedResource"

# Row 7 / 9 : native static handles (name/type/modifier interleaved)
$ws.Range("C7").Value = "g_outputTex"
$ws.Range("B9").Value = "cudaGraphicsResource*"
$ws.Range("B7").Value = "ComPtr<ID3D11Texture2D>"
$ws.Range("A7").Value = "static"

# Row 8
$ws.Range("A8").Value = "CudaRegisterD3D11Texture"

$ws.Range("A9").Value = "static"
$ws.Range("C9").Value = "g_cudaOut"

# Row 10
$ws.Range("A10").Value = "cudaGraphicsSubResourceGetMappedArray"

# Row 11
$ws.Range("B11").Value = "outArray"
$ws.Range("A11").Value = "void**"

# --- formatting ---------------------------------------------------------
# Apply left/top/wrap formatting only to the cells that actually hold
# content (the sheet otherwise has no populated cells outside these).
$populated = "A2","D2","F2","H2", `
             "A3","B3","C3", `
             "A4","B4","C4", `
             "A5","B5","C5","D5","E5","F5","G5","H5","I5","J5","K5", `
             "A6", `
             "A7","B7","C7", `
             "A8", `
             "A9","B9","C9", `
             "A10", `
             "A11","B11"
foreach ($addr in $populated) {
    $c = $ws.Range($addr)
    $c.HorizontalAlignment = -4131   # xlLeft
    $c.VerticalAlignment = -4160     # xlTop
    $c.WrapText = $true
}

$ws.Range("A2:K9").RowHeight = 37.5
$ws.Range("A10").RowHeight = 56.25

# --- selection / activation ---------------------------------------------
$ws.Range("C3").Select()
$ws.Activate()
